# Increase iron (column I) and IQ (column J) for male subjects.
# "male" appears in several case variants in the sex column (C): male, m, MALE.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $sex = $ws.Cells.Item($r, 3).Value()
    if ($sex -ceq "male" -or $sex -ceq "m" -or $sex -ceq "MALE") {
        $ironCell = $ws.Cells.Item($r, 9)
        $iqCell = $ws.Cells.Item($r, 10)

        $iron = $ironCell.Value()
        $iq = $iqCell.Value()

        $ironCell.Value = $iron + 2
        $iqCell.Value = $iq + 22
    }
}
